$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ofertas")

$ws.Range("A5").Value = "Empresa 3"
$ws.Range("B5").Value = "Articulo 57"
$ws.Range("C5").Value = 2000

$ws.Range("D5").Select()
